$d = $word.ActiveDocument
$tab = [char]9

# ---------------------------------------------------------------------------
# Edit 1: "...umber:    " + <w:tab/>  ->  "...umber:" + two literal spaces
#   (the paragraph already starts with a separate "N" run that stays intact)
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("umber:    " + $tab, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $start1 = $rng1.Start
    $end1 = $rng1.End

    # First 6 characters are "umber:" - touch (no-op bold toggle) to force a
    # run boundary right after the colon without altering its text/formatting.
    $partA = $d.Range($start1, $start1 + 6)
    $partA.Bold = 1
    $partA.Bold = 0

    # Remainder ("    " + tab) becomes two plain spaces in its own run.
    $partB = $d.Range($partA.End, $end1)
    $partB.Text = "  "
    $partB.Bold = 1
    $partB.Bold = 0
}

# ---------------------------------------------------------------------------
# Edit 2: "SOS Prepared By:" + <w:tab/>  ->  "SOS Prepared By" + ":  "
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("SOS Prepared By:" + $tab, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $start2 = $rng2.Start
    $end2 = $rng2.End

    # First 15 characters are "SOS Prepared By" - touch (no-op bold toggle)
    # to force a run boundary right before the colon.
    $partC = $d.Range($start2, $start2 + 15)
    $partC.Bold = 1
    $partC.Bold = 0

    # Remainder (":" + tab) becomes ":  " (colon + two spaces) in its own run.
    $partD = $d.Range($partC.End, $end2)
    $partD.Text = ":  "
    $partD.Bold = 1
    $partD.Bold = 0
}
